$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for rows 2-12 from 45212 to 45221 (date serial values)
$ws.Range("C2:C12").Value = 45221
